$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G3").Value = 0.8

$ws.Range("G14").Value = 0.7067669172932332

$ws.Range("F15").Value = 0.5729323308270677
$ws.Range("G15").Value = 0.6977443609022556

$ws.Range("G16").Value = 0.6962406015037594

$ws.Range("F17").Value = 0.5684210526315789
$ws.Range("G17").Value = 0.6947368421052631

$ws.Range("G19").Value = 0.8844444444444445

$ws.Range("F20").Value = 0.7999722222222223
$ws.Range("G20").Value = 0.8766666666666667

$ws.Range("F21").Value = 0.7933333333333333

$ws.Range("F22").Value = 0.8934239130434783

$ws.Range("F24").Value = 0.8673913043478261
$ws.Range("G24").Value = 0.9630434782608696

$ws.Range("F25").Value = 0.8695652173913043
$ws.Range("G25").Value = 0.9630434782608696
